$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "아핀 변환 (Affine Transformation)"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2024/06/28/Affine_Transformation.html"

$ws.Range("D28").Value = "통신 모델 및 관계 간단 정리"
$ws.Range("E28").Value = "https://ropiens.tistory.com/251"

$ws.Range("D51").Value = "[Oracle] 인덱스 조회 쿼리"
$ws.Range("E51").Value = "https://bskyvision.com/entry/Oracle-%EC%9D%B8%EB%8D%B1%EC%8A%A4-%EC%A1%B0%ED%9A%8C-%EC%BF%BC%EB%A6%AC"
